$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns: PREREQ_COURSES -> Prerequisites, COREQ_COURSES -> CoRequisites
$ws.Range("E1").Value = "Prerequisites"
$ws.Range("F1").Value = "CoRequisites"

# Update the selection to match the new selected range A1:I1
$ws.Range("A1:I1").Select()
